$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 196; $i -le 200; $i++) {
    $ws.Cells.Item($i, 1).Value = $i
}
